$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avichal")
$ws.Range("A8").Value = 45792
Write-Host "done"
